$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.317.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.587.66'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '506.53'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.68'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -6.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.595.11'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.61'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +7.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.347'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.89%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.044.38'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.336.73'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.60'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.31%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.596.14'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '347.21'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.36'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.22'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.703.07'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0845'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.42'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.31'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '154.37'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.55'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.75'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.00'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.45%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.848'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.46%  '
$ws.Range("B39").Value = 'SuiNetwork'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.842'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +18.30%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.50%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.82'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '297.03'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.71%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0999'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.02%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.618'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0561'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.73'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.88'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.69%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.09%  '
